$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMSProd")

# Update the provisioning identifiers for the LMS Prod setup config row.
$ws.Range("D2").Value = "fpadmin"
$ws.Range("A2").Value = "FPK12School83955"
$ws.Range("B2").Value = "FPK12Classroom26964"
$ws.Range("C2").Value = "FPK12Section94478"

# Move the active selection to the school/classroom/section block.
[void]$ws.Range("D2:F2").Select()
